# Generate Report for Handoff
# Adds two new entries (10c728f3-... and 993b042f-...) to the localization
# status workbook, on all three sheets (Overview, zh-cn, de-de), pushing the
# existing ".localization-config" row down by two rows.

$wb = $excel.ActiveWorkbook

# Blue hyperlink font color (RGB FF6495ED) expressed as the OLE/VBA BGR
# integer Excel's Font.Color expects.
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

function Style-AsDate($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

function Add-Hyperlink($ws, $addr, $url, $display) {
    $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $display)
    Style-AsHyperlink $ws.Range($addr)
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Rows.Item(4).Insert()
$ws1.Rows.Item(4).Insert()

# Drop every hyperlink on the sheet; they get rebuilt below in final order
# (row-insert does not shift existing hyperlink anchors, so the safest fix
# is to recreate them all from scratch).
$ws1.Range("A1").Hyperlinks.Delete()

$ws1.Range("A4").Value = "10c728f3-98e9-4b60-8915-d42e8cbab045.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = "993b042f-1602-41e2-a862-8e2b8a010545.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

Add-Hyperlink $ws1 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/415963f9-a38e-412d-8d9a-555ffcda4c7f.md" "415963f9-a38e-412d-8d9a-555ffcda4c7f.md"
Add-Hyperlink $ws1 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/777ff083-3218-4260-98b7-69434b99dba2.md" "777ff083-3218-4260-98b7-69434b99dba2.md"
Add-Hyperlink $ws1 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/10c728f3-98e9-4b60-8915-d42e8cbab045.md" "10c728f3-98e9-4b60-8915-d42e8cbab045.md"
Add-Hyperlink $ws1 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/993b042f-1602-41e2-a862-8e2b8a010545.md" "993b042f-1602-41e2-a862-8e2b8a010545.md"
Add-Hyperlink $ws1 "A6" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/.localization-config" ".localization-config"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(4).Insert()
$ws2.Rows.Item(4).Insert()

$ws2.Range("A1").Hyperlinks.Delete()

$ws2.Range("A4").Value = "10c728f3-98e9-4b60-8915-d42e8cbab045.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "10c728f3-98e9-4b60-8915-d42e8cbab045.52929bf5ce8ffb3eb1f66308197ca0ea0645777c.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-01-27 07:30:02"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Include"

$ws2.Range("A5").Value = "993b042f-1602-41e2-a862-8e2b8a010545.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = "993b042f-1602-41e2-a862-8e2b8a010545.09439e36c29e50089a4e6a252c3db8f424ca3470.zh-cn.xlf"
$ws2.Range("D5").Value = "2016-01-27 07:30:02"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Include"

$ws2.Range("D6").Value = "0001-01-01 00:00:00"
$ws2.Range("G6").Value = "0001-01-01 00:00:00"
$ws2.Range("H6").Value = "Ignored"

Add-Hyperlink $ws2 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/415963f9-a38e-412d-8d9a-555ffcda4c7f.md" "415963f9-a38e-412d-8d9a-555ffcda4c7f.md"
Add-Hyperlink $ws2 "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eaf100b68630fc304fdadb02bbfcd5cd0a223da5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/415963f9-a38e-412d-8d9a-555ffcda4c7f.cdb47c742a163eb9f6b00549f6726049e81613fa.zh-cn.xlf" "415963f9-a38e-412d-8d9a-555ffcda4c7f.cdb47c742a163eb9f6b00549f6726049e81613fa.zh-cn.xlf"
Add-Hyperlink $ws2 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/777ff083-3218-4260-98b7-69434b99dba2.md" "777ff083-3218-4260-98b7-69434b99dba2.md"
Add-Hyperlink $ws2 "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eaf100b68630fc304fdadb02bbfcd5cd0a223da5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/777ff083-3218-4260-98b7-69434b99dba2.c7b1442b54ced713484369667ec3d8c36a2aa057.zh-cn.xlf" "777ff083-3218-4260-98b7-69434b99dba2.c7b1442b54ced713484369667ec3d8c36a2aa057.zh-cn.xlf"
Add-Hyperlink $ws2 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/10c728f3-98e9-4b60-8915-d42e8cbab045.md" "10c728f3-98e9-4b60-8915-d42e8cbab045.md"
Add-Hyperlink $ws2 "C4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eaf100b68630fc304fdadb02bbfcd5cd0a223da5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/10c728f3-98e9-4b60-8915-d42e8cbab045.52929bf5ce8ffb3eb1f66308197ca0ea0645777c.zh-cn.xlf" "10c728f3-98e9-4b60-8915-d42e8cbab045.52929bf5ce8ffb3eb1f66308197ca0ea0645777c.zh-cn.xlf"
Add-Hyperlink $ws2 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/993b042f-1602-41e2-a862-8e2b8a010545.md" "993b042f-1602-41e2-a862-8e2b8a010545.md"
Add-Hyperlink $ws2 "C5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eaf100b68630fc304fdadb02bbfcd5cd0a223da5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/993b042f-1602-41e2-a862-8e2b8a010545.09439e36c29e50089a4e6a252c3db8f424ca3470.zh-cn.xlf" "993b042f-1602-41e2-a862-8e2b8a010545.09439e36c29e50089a4e6a252c3db8f424ca3470.zh-cn.xlf"
Add-Hyperlink $ws2 "A6" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/.localization-config" ".localization-config"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(4).Insert()
$ws3.Rows.Item(4).Insert()

$ws3.Range("A1").Hyperlinks.Delete()

$ws3.Range("A4").Value = "10c728f3-98e9-4b60-8915-d42e8cbab045.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "10c728f3-98e9-4b60-8915-d42e8cbab045.52929bf5ce8ffb3eb1f66308197ca0ea0645777c.de-de.xlf"
$ws3.Range("D4").Value = "2016-01-27 07:30:17"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Include"

$ws3.Range("A5").Value = "993b042f-1602-41e2-a862-8e2b8a010545.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = "993b042f-1602-41e2-a862-8e2b8a010545.09439e36c29e50089a4e6a252c3db8f424ca3470.de-de.xlf"
$ws3.Range("D5").Value = "2016-01-27 07:30:17"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Include"

$ws3.Range("D6").Value = "0001-01-01 00:00:00"
$ws3.Range("G6").Value = "0001-01-01 00:00:00"
$ws3.Range("H6").Value = "Ignored"

Add-Hyperlink $ws3 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/415963f9-a38e-412d-8d9a-555ffcda4c7f.md" "415963f9-a38e-412d-8d9a-555ffcda4c7f.md"
Add-Hyperlink $ws3 "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/989130b4f878b181b4569d68f27b05bb09d2f6ed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/415963f9-a38e-412d-8d9a-555ffcda4c7f.cdb47c742a163eb9f6b00549f6726049e81613fa.de-de.xlf" "415963f9-a38e-412d-8d9a-555ffcda4c7f.cdb47c742a163eb9f6b00549f6726049e81613fa.de-de.xlf"
Add-Hyperlink $ws3 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/777ff083-3218-4260-98b7-69434b99dba2.md" "777ff083-3218-4260-98b7-69434b99dba2.md"
Add-Hyperlink $ws3 "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/989130b4f878b181b4569d68f27b05bb09d2f6ed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/777ff083-3218-4260-98b7-69434b99dba2.c7b1442b54ced713484369667ec3d8c36a2aa057.de-de.xlf" "777ff083-3218-4260-98b7-69434b99dba2.c7b1442b54ced713484369667ec3d8c36a2aa057.de-de.xlf"
Add-Hyperlink $ws3 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/10c728f3-98e9-4b60-8915-d42e8cbab045.md" "10c728f3-98e9-4b60-8915-d42e8cbab045.md"
Add-Hyperlink $ws3 "C4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/989130b4f878b181b4569d68f27b05bb09d2f6ed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/10c728f3-98e9-4b60-8915-d42e8cbab045.52929bf5ce8ffb3eb1f66308197ca0ea0645777c.de-de.xlf" "10c728f3-98e9-4b60-8915-d42e8cbab045.52929bf5ce8ffb3eb1f66308197ca0ea0645777c.de-de.xlf"
Add-Hyperlink $ws3 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/e2e/993b042f-1602-41e2-a862-8e2b8a010545.md" "993b042f-1602-41e2-a862-8e2b8a010545.md"
Add-Hyperlink $ws3 "C5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/989130b4f878b181b4569d68f27b05bb09d2f6ed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/993b042f-1602-41e2-a862-8e2b8a010545.09439e36c29e50089a4e6a252c3db8f424ca3470.de-de.xlf" "993b042f-1602-41e2-a862-8e2b8a010545.09439e36c29e50089a4e6a252c3db8f424ca3470.de-de.xlf"
Add-Hyperlink $ws3 "A6" "https://github.com/OpenLocalizationTest/oltest/blob/ebe54c77a3c6664e7f01ea1ca0354a5af6191c8d/.localization-config" ".localization-config"
